$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (F column) for rows 3-6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1879
$ws1.Range("F4").Value = 816
$ws1.Range("F5").Value = 698
$ws1.Range("F6").Value = 235

# Sheet "全部类型" - update "想去人数" (F column) for rows 3,5,6,7
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1879
$ws4.Range("F5").Value = 816
$ws4.Range("F6").Value = 698
$ws4.Range("F7").Value = 235
